# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet (which mirrors the same data).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 483
    $ws.Range("F3").Value = 70
    $ws.Range("F4").Value = 43
}
